$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Date: updated publication date
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: false -> true
# NB: setting Range.Value = "true" directly gets auto-coerced by Excel into the
# boolean TRUE (same quirk real Excel has), which would store a <c t="b"> cell
# instead of the literal text "true". Route the literal string through a
# formula + PasteSpecial(values) round-trip so it lands as plain text.
$helper = $ws.Range("Z1")
$helper.Formula = '=TRIM("true ")'
$target = $ws.Range("B17")
$helper.Copy()
$target.PasteSpecial(-4163)
$helper.Clear()
$excel.CutCopyMode = $false
